$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversion del dia" note with today's rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 14.49 = 60019.71 pesos`n✅ 60019.71 pesos = 14.5 = 981.98 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- tasas: refresh the raw exchange-rate inputs ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 68.994
$ws2.Range("O10").Value = 4141
$ws2.Range("N12").Value = 4137.91
$ws2.Range("O12").Value = 67.7
